$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- New "Fractions" section header (row 45/46) ---
$ws1.Range("G45").Value = "Fractions"

$ws1.Range("G46").Value = "ZEROS"
$ws1.Range("H46").Value = "REPEATED"
$ws1.Range("I46").Value = "MEDIUM"
$ws1.Range("J46").Value = "LARGE"
$hdr = $ws1.Range("G46:J46")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108

# --- Updated raw counts (rows 51-52) ---
$ws1.Range("D51").Value = 1
$ws1.Range("E51").Value = 2047
$ws1.Range("B52").Value = 0
$ws1.Range("D52").Value = 1251
$ws1.Range("E52").Value = 782
$ws1.Range("F52").Value = 2033

# --- New fraction-of-total formulas (G:J, rows 47-52) ---
$ws1.Range("G47").Formula = "=B47/F47"
$ws1.Range("H47").Formula = "=C47/F47"
$ws1.Range("I47").Formula = "=D47/F47"
$ws1.Range("J47").Formula = "=E47/F47"

$ws1.Range("G48:G52").Formula = "=B48/F48"
$ws1.Range("H48:H52").Formula = "=C48/F48"
$ws1.Range("I48:I52").Formula = "=D48/F48"
$ws1.Range("J48:J52").Formula = "=E48/F48"

$pct = $ws1.Range("G47:J52")
$pct.HorizontalAlignment = -4152
$pct.Font.Bold = $false
$pct.NumberFormat = "0.00%"
